$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# row 15
$ws.Range("H15").Value = 1942.3636
$ws.Range("I15").Value = 1942.3636
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 5827.0908
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -5658.0908

# row 28
$ws.Range("H28").Value = 1477.3125
$ws.Range("I28").Value = 979.8333
$ws.Range("J28").Value = 2969.75
$ws.Range("K28").Value = 979.8333
$ws.Range("L28").Value = 2969.75
$ws.Range("M28").Value = -494.8333

# row 32
$ws.Range("H32").Value = 8183
$ws.Range("I32").Value = 8386.5
$ws.Range("J32").Value = 7776
$ws.Range("K32").Value = 8386.5
$ws.Range("L32").Value = 7776
$ws.Range("M32").Value = -8060.5
$ws.Range("N32").Value = -8428

# row 37
$ws.Range("H37").Value = 400
$ws.Range("I37").Value = 100
$ws.Range("J37").Value = 1000
$ws.Range("K37").Value = 300
$ws.Range("L37").Value = 3000
$ws.Range("M37").Value = -174

# row 43
$ws.Range("H43").Value = 5637.75
$ws.Range("I43").Value = 5683.6665
$ws.Range("J43").Value = 5500
$ws.Range("K43").Value = 5683.6665
$ws.Range("L43").Value = 5500
$ws.Range("M43").Value = -5614.6665

# row 74
$ws.Range("H74").Value = 4062
$ws.Range("I74").Value = 4062
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 4062
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -3126

# row 77
$ws.Range("H77").Value = 4062
$ws.Range("I77").Value = 4062
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 20310
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -15630

# row 86
$ws.Range("H86").Value = 4000.6
$ws.Range("I86").Value = 4500
$ws.Range("J86").Value = 3875.75
$ws.Range("K86").Value = 4500
$ws.Range("L86").Value = 3875.75
$ws.Range("M86").Value = -3377

# row 89
$ws.Range("H89").Value = 4000.6
$ws.Range("I89").Value = 4500
$ws.Range("J89").Value = 3875.75
$ws.Range("K89").Value = 22500
$ws.Range("L89").Value = 19378.75
$ws.Range("M89").Value = -16884

# row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()

# row 116
$ws.Range("H116").Value = 6333
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 6333
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 6333
$ws.Range("M116").ClearContents()
$ws.Range("N116").Value = -13217

# row 132
$ws.Range("H132").Value = 4912.1787
$ws.Range("I132").Value = 4960.926
$ws.Range("J132").Value = 3596
$ws.Range("K132").Value = 14882.778
$ws.Range("L132").Value = 10788
$ws.Range("M132").Value = -12352.778
$ws.Range("N132").Value = -15848

# row 137
$ws.Range("H137").Value = 2106.4614
$ws.Range("I137").Value = 2106.4614
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 6319.3842
$ws.Range("L137").Value = 0
$ws.Range("M137").Value = -3769.3842

$ws = $wb.Worksheets.Item("ARM")
# row 36
$ws.Range("H36").Value = 15000
$ws.Range("I36").Value = 15000
$ws.Range("J36").Value = 0
$ws.Range("K36").Value = 15000
$ws.Range("L36").Value = 0
$ws.Range("M36").Value = -14654

# row 74
$ws.Range("H74").Value = 1177.6
$ws.Range("I74").Value = 982.4286
$ws.Range("J74").Value = 1633
$ws.Range("K74").Value = 982.4286
$ws.Range("L74").Value = 1633
$ws.Range("M74").Value = -108.4286
$ws.Range("N74").Value = -3381

# row 77
$ws.Range("H77").Value = 1177.6
$ws.Range("I77").Value = 982.4286
$ws.Range("J77").Value = 1633
$ws.Range("K77").Value = 4912.143
$ws.Range("L77").Value = 8165
$ws.Range("M77").Value = -544.143
$ws.Range("N77").Value = -16901

# row 132
$ws.Range("H132").Value = 2320.111
$ws.Range("I132").Value = 2320.111
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 6960.333
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4430.333

$ws = $wb.Worksheets.Item("BSM")
# row 7
$ws.Range("H7").Value = 0
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 0
$ws.Range("N7").ClearContents()

# row 69
$ws.Range("H69").Value = 30000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 30000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 30000
$ws.Range("N69").Value = -31622

# row 72
$ws.Range("H72").Value = 30000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 30000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 90000
$ws.Range("N72").Value = -98112

$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 1321
$ws.Range("I31").Value = 1321
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 1321
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -1026
$ws.Range("N31").ClearContents()

# row 34
$ws.Range("H34").Value = 1321
$ws.Range("I34").Value = 1321
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 1321
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -1119
$ws.Range("N34").ClearContents()

# row 132
$ws.Range("H132").Value = 1846.4615
$ws.Range("I132").Value = 1592.9166
$ws.Range("J132").Value = 4889
$ws.Range("K132").Value = 4778.7498
$ws.Range("L132").Value = 14667
$ws.Range("M132").Value = -2248.7498
$ws.Range("N132").Value = -19727

$ws = $wb.Worksheets.Item("CUL")
# row 25
$ws.Range("H25").Value = 693.75
$ws.Range("I25").Value = 693.75
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 2081.25
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = -1912.25
$ws.Range("N25").ClearContents()

# row 30
$ws.Range("H30").Value = 693.75
$ws.Range("I30").Value = 693.75
$ws.Range("J30").Value = 0
$ws.Range("K30").Value = 2081.25
$ws.Range("L30").Value = 0
$ws.Range("M30").Value = -1979.25
$ws.Range("N30").ClearContents()

# row 106
$ws.Range("H106").Value = 19165.111
$ws.Range("I106").Value = 9990
$ws.Range("J106").Value = 20312
$ws.Range("K106").Value = 29970
$ws.Range("L106").Value = 60936
$ws.Range("M106").Value = -29024
$ws.Range("N106").Value = -62828

# row 111
$ws.Range("H111").Value = 5166.3335
$ws.Range("I111").Value = 5166.3335
$ws.Range("J111").Value = 0
$ws.Range("K111").Value = 15499.0005
$ws.Range("L111").Value = 0
$ws.Range("M111").Value = -12432.0005

$ws = $wb.Worksheets.Item("GSM")
# row 57
$ws.Range("H57").Value = 25000
$ws.Range("I57").Value = 25000
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 25000
$ws.Range("L57").Value = 0
$ws.Range("M57").Value = -24180

# row 58
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()

# row 88
$ws.Range("H88").Value = 50000
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 50000
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 50000
$ws.Range("N88").Value = -50902

# row 91
$ws.Range("H91").Value = 50000
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 50000
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 50000
$ws.Range("N91").Value = -53120

# row 102
$ws.Range("H102").Value = 1366.75
$ws.Range("I102").Value = 1366.75
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 1366.75
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = 255.25

# row 113
$ws.Range("H113").Value = 0
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()

# row 128
$ws.Range("H128").Value = 47196.6
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 47196.6
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 47196.6
$ws.Range("M128").ClearContents()
$ws.Range("N128").Value = -57156.6

# row 132
$ws.Range("H132").Value = 2129
$ws.Range("I132").Value = 1968.8462
$ws.Range("J132").Value = 2649.5
$ws.Range("K132").Value = 5906.5386
$ws.Range("L132").Value = 7948.5
$ws.Range("M132").Value = -3376.5386
$ws.Range("N132").Value = -13008.5

$ws = $wb.Worksheets.Item("LTW")
# row 7
$ws.Range("H7").Value = 6695.2
$ws.Range("I7").Value = 6695.2
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 6695.2
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -6583.2

# row 40
$ws.Range("H40").Value = 1962.3334
$ws.Range("I40").Value = 1962.3334
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1962.3334
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1826.3334

# row 46
$ws.Range("H46").Value = 1364.2858
$ws.Range("I46").Value = 1314
$ws.Range("J46").Value = 1490
$ws.Range("K46").Value = 1314
$ws.Range("L46").Value = 1490
$ws.Range("M46").Value = -1126
$ws.Range("N46").Value = -1866

# row 126
$ws.Range("H126").Value = 6695.2
$ws.Range("I126").Value = 6695.2
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 20085.6
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -17615.6

# row 136
$ws.Range("H136").Value = 1616
$ws.Range("I136").Value = 1616
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 4848
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -2298

$ws = $wb.Worksheets.Item("WVR")
# row 52
$ws.Range("H52").Value = 22268.6
$ws.Range("I52").Value = 18135.75
$ws.Range("J52").Value = 38800
$ws.Range("K52").Value = 18135.75
$ws.Range("L52").Value = 38800
$ws.Range("M52").Value = -17909.75

# row 58
$ws.Range("H58").Value = 50497.5
$ws.Range("I58").Value = 6995
$ws.Range("J58").Value = 94000
$ws.Range("K58").Value = 6995
$ws.Range("L58").Value = 94000
$ws.Range("M58").Value = -6687

# row 62
$ws.Range("H62").Value = 8766.625
$ws.Range("I62").Value = 10050
$ws.Range("J62").Value = 8338.833000000001
$ws.Range("K62").Value = 10050
$ws.Range("L62").Value = 8338.833000000001
$ws.Range("M62").Value = -9426
$ws.Range("N62").Value = -9586.833000000001

# row 65
$ws.Range("H65").Value = 8766.625
$ws.Range("I65").Value = 10050
$ws.Range("J65").Value = 8338.833000000001
$ws.Range("K65").Value = 50250
$ws.Range("L65").Value = 41694.165
$ws.Range("M65").Value = -47130
$ws.Range("N65").Value = -47934.165

# row 99
$ws.Range("H99").Value = 50475
$ws.Range("I99").Value = 50475
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 50475
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -47480
$ws.Range("N99").ClearContents()

# row 132
$ws.Range("H132").Value = 3479.4
$ws.Range("I132").Value = 3549.3572
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 10648.0716
$ws.Range("L132").Value = 7500
$ws.Range("M132").Value = -8118.071599999999

# row 136
$ws.Range("H136").Value = 7533.185
$ws.Range("I136").Value = 4854.1816
$ws.Range("J136").Value = 9375
$ws.Range("K136").Value = 14562.5448
$ws.Range("L136").Value = 28125
$ws.Range("M136").Value = -12012.5448
